$d = $word.ActiveDocument

# --- Locate the two anchor paragraphs reliably by content, skipping table paragraphs ---
$count = $d.Paragraphs.Count
$freeFormIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.Information(12) -eq $false) -and ($p.Range.Text -like "*<<freeFormRecitalText>>*")) {
        $freeFormIndex = $i
        break
    }
}
if ($freeFormIndex -eq -1) {
    throw "Could not locate the <<freeFormRecitalText>> paragraph"
}

$recitalEmptyIndex = $freeFormIndex - 2

# --- Edit 1: turn the blank paragraph before "THE COURT RECORDS THAT:" into the
#             new Docmosis conditional-section-start tag for the recital text ---
$emptyPara = $d.Paragraphs.Item($recitalEmptyIndex)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>&lt;&lt;cs</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>_</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>freeFormRecitalText</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t xml:space="preserve"> != null}</w:t></w:r><w:r><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>&gt;&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$emptyPara.Range.InsertXML($xml1)

# --- Edit 2: insert a new paragraph holding the matching conditional-section-end
#             tag right after the <<freeFormRecitalText>> paragraph ---
$freeFormPara = $d.Paragraphs.Item($freeFormIndex)
$freeFormPara.Range.InsertParagraphAfter()
$esPara = $d.Paragraphs.Item($freeFormIndex + 1)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="2"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-US"/><w14:ligatures w14:val="standardContextual"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="2"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-US"/><w14:ligatures w14:val="standardContextual"/></w:rPr><w:lastRenderedPageBreak/><w:t>&lt;&lt;es_&gt;&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$esPara.Range.InsertXML($xml2)

Write-Host "Edits applied."
